# The deck's single embedded "live" theme (ppt/theme/theme2.xml, referenced
# by the slide master) currently carries the "Integral" color scheme.
# The target state swaps it to the "Office Theme" color scheme
# (font scheme / format scheme are already identical between the two
# theme parts, so only the 12 scheme colors need to change).
#
# PowerPoint COM RGB values are packed little-endian as 0x00BBGGRR, so
# each target color is rebuilt from its R/G/B bytes with bit shifts
# instead of hard-coding opaque decimal constants.

function New-RgbValue($r, $g, $b) {
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$p = $ppt.ActivePresentation
$theme = $p.Designs.Item(1).SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$colors.Colors(1).RGB  = New-RgbValue 0x00 0x00 0x00   # dk1      000000
$colors.Colors(2).RGB  = New-RgbValue 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colors.Colors(3).RGB  = New-RgbValue 0x44 0x54 0x6A   # dk2      44546A
$colors.Colors(4).RGB  = New-RgbValue 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colors.Colors(5).RGB  = New-RgbValue 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colors.Colors(6).RGB  = New-RgbValue 0xED 0x7D 0x31   # accent2  ED7D31
$colors.Colors(7).RGB  = New-RgbValue 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colors.Colors(8).RGB  = New-RgbValue 0xFF 0xC0 0x00   # accent4  FFC000
$colors.Colors(9).RGB  = New-RgbValue 0x44 0x72 0xC4   # accent5  4472C4
$colors.Colors(10).RGB = New-RgbValue 0x70 0xAD 0x47   # accent6  70AD47
$colors.Colors(11).RGB = New-RgbValue 0x05 0x63 0xC1   # hlink    0563C1
$colors.Colors(12).RGB = New-RgbValue 0x95 0x4F 0x72   # folHlink 954F72
